# add try and exceptions
# Update totalReports / numDistinctUsers counts and lastReportedAt timestamps
# for several rows in the IP-blocked report after re-running the lookup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    # Row 5 - 198.211.120.36
    $ws.Range("K5").Value = 772
    $ws.Range("L5").Value = 350
    $ws.Range("M5").Value = "2023-08-27T15:13:43+00:00"

    # Row 10 - 201.236.101.194
    $ws.Range("M10").Value = "2023-08-27T15:08:47+00:00"

    # Row 11 - 79.124.62.82
    $ws.Range("K11").Value = 14278
    $ws.Range("M11").Value = "2023-08-27T15:04:45+00:00"

    # Row 15 - 185.126.80.63
    $ws.Range("K15").Value = 198
    $ws.Range("L15").Value = 123
    $ws.Range("M15").Value = "2023-08-27T15:00:03+00:00"

    # Row 16 - 195.33.237.83
    $ws.Range("K16").Value = 2994

    # Row 17 - 51.210.254.243
    $ws.Range("K17").Value = 1493
    $ws.Range("L17").Value = 519
    $ws.Range("M17").Value = "2023-08-27T15:10:18+00:00"

    # Row 19 - 143.42.164.127
    $ws.Range("K19").Value = 941
    $ws.Range("M19").Value = "2023-08-27T15:03:37+00:00"

    # Row 21 - 159.65.194.58
    $ws.Range("K21").Value = 3039
    $ws.Range("M21").Value = "2023-08-27T15:13:21+00:00"

    # Row 22 - 140.99.4.7
    $ws.Range("K22").Value = 1066
    $ws.Range("L22").Value = 429
    $ws.Range("M22").Value = "2023-08-27T15:05:31+00:00"
}
catch {
    Write-Host "Error while updating report: $_"
}
